$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = "P. point"
$ws.Range("C8").Value = 82
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "Medium point (up to 6 mtr.)"
$ws.Range("F8").Value = 472
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "38704.00"
# Row 9
$ws.Range("C9").Value = 64
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "4"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "Long point  (up to 10 mtr.)"
$ws.Range("F9").Value = 662
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "42368.00"
# Row 10
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = "Each"
$ws.Range("C10").Value = 20
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "4.0"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "P & F ISI marked (IS :3854) 16 amp. flush type non modular switch CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including cutting hole in tile and making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure `"A`" attached with this BSR ."
$ws.Range("F10").Value = 50
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "1000.00"
# Row 11
$ws.Range("A11").NumberFormat = "@"
$ws.Range("A11").Value = "Each"
$ws.Range("C11").Value = 56
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.0"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "Providing & Fixing of  3/6 pin 16 amp flush type non modular socket  made out from Industrial grade Polycarbonate or fire resistant ABS material, brass terminal with Porcelain based back cover & captive screws including cutting hole in tile and making connection testing etc. as required.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure `"A`" attached with this BSR ."
$ws.Range("F11").Value = 78
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "4368.00"
# Row 12
$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = ""
$ws.Range("C12").Value = 48
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "11.0"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "S&F following sizes (dia.) of ISI marked virgin material MMS ( IS:9537 P - III ) PVC conduit along with  ISI marked (IS:3419-1988) accessories as required  in  recess  including  cutting the wall, covering conduit and making good the same as required. For additional technical parameters of product / work refer Annexure 'A' attached with this BSR"
$ws.Range("F12").Value = 0
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "0.00"
# Row 13
$ws.Range("A13").NumberFormat = "@"
$ws.Range("A13").Value = "R. mtr."
$ws.Range("C13").Value = 27
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "17"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "25 mm"
$ws.Range("F13").Value = 56
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "1512.00"
# Row 14
$ws.Range("A14").NumberFormat = "@"
$ws.Range("A14").Value = "Set"
$ws.Range("C14").Value = 48
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "13.0"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "Plate Earthing  as per IS:3043 with Hot dipped G.I. Earth plate of size 600mm x 600mm x 6.0mm by embodying  3 to 4 mtr. below the ground level with 20  mm dia. G.I. 'B' class watering Pipe ,including all accessories like nut, bolts, reducer, nipple, wire meshed funnel, and Heavy duty weather proof poly-propylene earth pit chamber with lockable Jam free lid suitable for safe working load 5000 Kg or more of size Top Dia. 225 to 260 mm, Bottom Dia 300 to 350 mm. and Height  250 to 300 mm. and embodying the pipe  complete with alternate layers salt and coke/ charcoal, testing of earth resistance for value of 5 ohms or less  as required & must record by engineer in charge during site visit and ensure to enter in measurment book.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure `"A`" attached with this BSR .   "
$ws.Range("F14").Value = 5733
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "275184.00"
# Row 15
$ws.Range("A15").NumberFormat = "@"
$ws.Range("A15").Value = ""
$ws.Range("C15").Value = 13
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "16.0"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "Providing & Fixing of IP20 SMD Mid Power LED batten type integrated light fixture made from Powder coated Extruded aluminium  housing with in built driver  , System lumen efficacy ≥ 110 lm/Watt output, internal surge protection of 2.5 KV with Short & Open circuit protection ,THD < 10% , P. F.≥0.95, CRI >80 , life time of minimum  50000 Burning Hours with , 70% of intial Lumen maintaned till life ends  , CCT 3000°K / 4000°K  / 5700°K /6000°K/6500°K (As per ANSI Bin) , Maximum power consumption should not more than the specified rating and Fixture shall be of  BIS standard and  trade mark certificate ( T.C.). Manufactures Word Mark/ Name Engraved/ Embossing/ Screen printing on housing. OEM must have its own in house NABL lab setup for all testing facilities for LED fixtures. (LM79 & LM80) certificate / Report from OEM shall be submitted.  All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure `"A`" attached with this BSR ."
$ws.Range("F15").Value = 0
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "0.00"
# Row 16
$ws.Range("A16").NumberFormat = "@"
$ws.Range("A16").Value = ""
$ws.Range("C16").Value = 9
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "17.0"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "Providing & Fixing of 240/415 V AC MCB with positive isolation of 10 kA breaking capacity (B/ C/D tripping characteristic as per type of load and  site requirement) 4 KV impulse withstand voltage, ISI marked IS 8828(1996) / conforming to IEC 60898-1 2002, IEC 60947-2, low watt losses, trip free mechanisum , energy limiting of  class 3 as per IEC,  minimum phase termination capacity of 35sq.mm. , conductor line load reversibility , IP 20 contact protection and fitted in  existing distribution board/sheets, minimum electrical operation 20,000 upto 20 A rating and 10,000 upto 63 A, 5000 for 80 A & above rating  including making connections, testing etc. as required. OEM shall have submit  NABL / CPRI / ERDA accrediated   lab type test reports  & All as per pre approved by Engineer in charge. For additional technical parameters of product / work refer Annexure 'A' attached with this BSR"
$ws.Range("F16").Value = 0
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "0.00"
# Row 17
$ws.Range("C17").Value = 87
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "31"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "Double pole MCB(With B/C curve tripping Characteristics)"
# Row 18
$ws.Range("C18").Value = 37
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "32"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = " 50/63 A rating"
$ws.Range("F18").Value = 900
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "33300.00"
# Row 19
$ws.Range("A19").NumberFormat = "@"
$ws.Range("A19").Value = "Each"
$ws.Range("C19").Value = 13
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "35"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "8 Way (8+2)"
$ws.Range("F19").Value = 2184
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "28392.00"
# Row 20
$ws.Range("A20").NumberFormat = "@"
$ws.Range("A20").Value = ""
$ws.Range("C20").Value = 38
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "36"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "Total"
# Row 21
$ws.Range("C21").Value = 35
# Row 23
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "424828.00"
$ws.Range("H23").NumberFormat = "@"
$ws.Range("H23").Value = "424828.00"
# Row 25
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "424828.00"
$ws.Range("H25").NumberFormat = "@"
$ws.Range("H25").Value = "424828.00"
